$wb = $excel.ActiveWorkbook

# Add the new worksheet "SWOPIT_NONE_TRUE_250" after the last existing sheet
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "SWOPIT_NONE_TRUE_250"

$ws.Range("A1").Value = "ESTIMATION RESULTS for SWOPIT_NONE_TRUE_250"
$ws.Range("I2").Value = "Hours"
$ws.Range("J2").Value = "Mins"
$ws.Range("K2").Value = "Secs"
$ws.Range("A3").Value = "Converged: "
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "Not Converged: "
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "Startiter:"
$ws.Range("F3").Value = 1
$ws.Range("H3").Value = "Runtime:"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4
$ws.Range("A4").Value = "PARAMETERS SWOPIT_NONE_TRUE_250"
$ws.Range("A5").Value = "NAMES"
$ws.Range("B5").Formula = "'TRUE"
$ws.Range("C5").Value = "mean"
$ws.Range("D5").Value = "mean ci_low"
$ws.Range("E5").Value = "mean ci_high"
$ws.Range("F5").Value = "mean se"
$ws.Range("G5").Value = "real se"
$ws.Range("H5").Value = "real2mean se"
$ws.Range("I5").Value = "real2median se"
$ws.Range("J5").Value = "rmse"
$ws.Range("K5").Value = "coverage"
$ws.Range("A6").Value = "g1"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 2.308942002382028
$ws.Range("D6").Value = 1.0296123777843826
$ws.Range("E6").Value = 3.5882716269796733
$ws.Range("F6").Value = 0.65273119031208449
$ws.Range("G6").Value = 0.19685624139616253
$ws.Range("H6").Value = 0.3015885318764091
$ws.Range("I6").Value = 0.3015885318764091
$ws.Range("J6").Value = 0.36632982490160582
$ws.Range("K6").Value = 1
$ws.Range("A7").Value = "reg_cutoff"
$ws.Range("B7").Value = 0.20000000000000001
$ws.Range("C7").Value = -0.023488126569383788
$ws.Range("D7").Value = -0.60106948767062574
$ws.Range("E7").Value = 0.55409323453185833
$ws.Range("F7").Value = 0.29468978290271164
$ws.Range("G7").Value = 0.027866256961893661
$ws.Range("H7").Value = 0.09456132712647651
$ws.Range("I7").Value = 0.09456132712647651
$ws.Range("J7").Value = 0.22521871812653407
$ws.Range("K7").Value = 1
$ws.Range("A8").Value = "b2"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 2.1204793560804007
$ws.Range("D8").Value = 1.098956938519267
$ws.Range("E8").Value = 3.1420017736415344
$ws.Range("F8").Value = 0.52119448398989587
$ws.Range("G8").Value = 0.13485424150793313
$ws.Range("H8").Value = 0.25874073047662466
$ws.Range("I8").Value = 0.25874073047662466
$ws.Range("J8").Value = 0.18083401697199541
$ws.Range("K8").Value = 1
$ws.Range("A9").Value = "b3"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0.99430707012951403
$ws.Range("D9").Value = 0.47123197774861564
$ws.Range("E9").Value = 1.5173821625104127
$ws.Range("F9").Value = 0.26687995111484097
$ws.Range("G9").Value = 0.039671795032221784
$ws.Range("H9").Value = 0.14865033835063404
$ws.Range("I9").Value = 0.14865033835063404
$ws.Range("J9").Value = 0.040078183237129
$ws.Range("K9").Value = 1
$ws.Range("A10").Value = "out1_cutoff1"
$ws.Range("B10").Value = -3.8300000000000001
$ws.Range("C10").Value = -3.8014574741731479
$ws.Range("D10").Value = -5.8196150135731299
$ws.Range("E10").Value = -1.7832999347731666
$ws.Range("F10").Value = 1.0296911347958184
$ws.Range("G10").Value = 0.11068148337835182
$ws.Range("H10").Value = 0.10748998378071818
$ws.Range("I10").Value = 0.10748998378071818
$ws.Range("J10").Value = 0.11430252203433167
$ws.Range("K10").Value = 1
$ws.Range("A11").Value = "out1_cutoff2"
$ws.Range("B11").Value = 3.7599999999999998
$ws.Range("C11").Value = 3.957347871891117
$ws.Range("D11").Value = 1.9065762250737548
$ws.Range("E11").Value = 6.0081195187084786
$ws.Range("F11").Value = 1.0463312912857521
$ws.Range("G11").Value = 0.30405860814442259
$ws.Range("H11").Value = 0.29059496803425378
$ws.Range("I11").Value = 0.29059496803425378
$ws.Range("J11").Value = 0.36248837185029292
$ws.Range("K11").Value = 1
$ws.Range("A12").Value = "b4"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 1.2972650094676323
$ws.Range("D12").Value = 0.54349234713706063
$ws.Range("E12").Value = 2.0510376717982037
$ws.Range("F12").Value = 0.38458495578298091
$ws.Range("G12").Value = 0.26412664696620752
$ws.Range("H12").Value = 0.6867836169734437
$ws.Range("I12").Value = 0.6867836169734437
$ws.Range("J12").Value = 0.39765483964287807
$ws.Range("K12").Value = 1
$ws.Range("A13").Value = "b5"
$ws.Range("B13").Value = -2
$ws.Range("C13").Value = -2.4891349731677042
$ws.Range("D13").Value = -3.8915905586309063
$ws.Range("E13").Value = -1.086679387704502
$ws.Range("F13").Value = 0.71555171244246996
$ws.Range("G13").Value = 0.35592377481144943
$ws.Range("H13").Value = 0.49741167356938659
$ws.Range("I13").Value = 0.49741167356938659
$ws.Range("J13").Value = 0.60492541313107517
$ws.Range("K13").Value = 1
$ws.Range("A14").Value = "out2_cutoff1"
$ws.Range("B14").Value = -3.9700000000000002
$ws.Range("C14").Value = -5.0050694401860794
$ws.Range("D14").Value = -7.792360258666827
$ws.Range("E14").Value = -2.2177786217053326
$ws.Range("F14").Value = 1.4221132839514103
$ws.Range("G14").Value = 0.86565680706225656
$ws.Range("H14").Value = 0.60871156808055926
$ws.Range("I14").Value = 0.60871156808055926
$ws.Range("J14").Value = 1.3493444532884642
$ws.Range("K14").Value = 1
$ws.Range("A15").Value = "out2_cutoff2"
$ws.Range("B15").Value = 3.9700000000000002
$ws.Range("C15").Value = 4.8973168503969733
$ws.Range("D15").Value = 2.000185461677221
$ws.Range("E15").Value = 7.7944482391167256
$ws.Range("F15").Value = 1.4781554210036285
$ws.Range("G15").Value = 0.88345568745685199
$ws.Range("H15").Value = 0.59767442239396507
$ws.Range("I15").Value = 0.59767442239396507
$ws.Range("J15").Value = 1.2807851079435695
$ws.Range("K15").Value = 1
$ws.Range("A18").Value = "PROBABILITIES SWOPIT_NONE_TRUE_250"
$ws.Range("A19").Value = "Choice/Reg"
$ws.Range("B19").Formula = "'TRUE"
$ws.Range("C19").Value = "mean"
$ws.Range("D19").Value = "mean ci_low"
$ws.Range("E19").Value = "mean ci_high"
$ws.Range("F19").Value = "mean se"
$ws.Range("G19").Value = "real se"
$ws.Range("H19").Value = "real2mean se"
$ws.Range("I19").Value = "real2median se"
$ws.Range("J19").Value = "rmse"
$ws.Range("K19").Value = "coverage"
$ws.Range("A20").Value = "Choice 1"
$ws.Range("B20").Value = 0.34081173787569552
$ws.Range("C20").Value = 0.42155510504461091
$ws.Range("D20").Value = 0.16425103686925702
$ws.Range("E20").Value = 0.67885917321996481
$ws.Range("F20").Value = 0.13127999810452418
$ws.Range("G20").Value = 0.019674932229182509
$ws.Range("H20").Value = 0.14986999172194893
$ws.Range("I20").Value = 0.14986999172194893
$ws.Range("J20").Value = 0.084628847363468498
$ws.Range("K20").Value = 1
$ws.Range("A21").Value = "Choice 2"
$ws.Range("B21").Value = 0.33322877508788151
$ws.Range("C21").Value = 0.4273943182615455
$ws.Range("D21").Value = 0.1097162034830334
$ws.Range("E21").Value = 0.74507243304005766
$ws.Range("F21").Value = 0.1620836491304517
$ws.Range("G21").Value = 0.0094599176904864267
$ws.Range("H21").Value = 0.058364417023167393
$ws.Range("I21").Value = 0.058364417023167393
$ws.Range("J21").Value = 0.10008465497157898
$ws.Range("K21").Value = 1
$ws.Range("A22").Value = "Choice 3"
$ws.Range("B22").Value = 0.32595948703642297
$ws.Range("C22").Value = 0.15105057669384356
$ws.Range("D22").Value = -0.092303511215936523
$ws.Range("E22").Value = 0.39440466460362367
$ws.Range("F22").Value = 0.12416253044919506
$ws.Range("G22").Value = 0.0016668143207748958
$ws.Range("H22").Value = 0.013424455145563615
$ws.Range("I22").Value = 0.013424455145563615
$ws.Range("J22").Value = 0.18467412462890859
$ws.Range("K22").Value = 1
$ws.Range("A23").Value = "Reg 1"
$ws.Range("B23").Value = 0.49139299470263359
$ws.Range("C23").Value = 0.39014271432940251
$ws.Range("D23").Value = 0.15827428145937553
$ws.Range("E23").Value = 0.6220111471994294
$ws.Range("F23").Value = 0.11830239468631856
$ws.Range("G23").Value = 0.0017787522764088418
$ws.Range("H23").Value = 0.015035640496756156
$ws.Range("I23").Value = 0.015035640496756156
$ws.Range("J23").Value = 0.1030272487236889
$ws.Range("K23").Value = 1
$ws.Range("A24").Value = "Reg 2"
$ws.Range("B24").Value = 0.50860700529736635
$ws.Range("C24").Value = 0.60985728567059749
$ws.Range("D24").Value = 0.37798885296105333
$ws.Range("E24").Value = 0.84172571838014165
$ws.Range("F24").Value = 0.11830239460443806
$ws.Range("G24").Value = 0.0017787523440540851
$ws.Range("H24").Value = 0.0150356410789622
$ws.Range("I24").Value = 0.0150356410789622
$ws.Range("J24").Value = 0.10302724872368894
$ws.Range("K24").Value = 1
$ws.Range("A27").Value = "MARGINAL EFFECT SWOPIT_NONE_TRUE_250"
$ws.Range("A28").Value = "NAMES"
$ws.Range("B28").Formula = "'TRUE"
$ws.Range("C28").Value = "mean"
$ws.Range("D28").Value = "mean ci_low"
$ws.Range("E28").Value = "mean ci_high"
$ws.Range("F28").Value = "mean se"
$ws.Range("G28").Value = "real se"
$ws.Range("H28").Value = "real2mean se"
$ws.Range("I28").Value = "real2median se"
$ws.Range("J28").Value = "rmse"
$ws.Range("K28").Value = "coverage"
$ws.Range("A29").Value = "X1 on 1"
$ws.Range("B29").Value = 0.53452888450843739
$ws.Range("C29").Value = 0.61165685329914998
$ws.Range("D29").Value = 0.17582289468608406
$ws.Range("E29").Value = 1.0474908119122157
$ws.Range("F29").Value = 0.22236835066912897
$ws.Range("G29").Value = 0.01393725596215974
$ws.Range("H29").Value = 0.062676437182814557
$ws.Range("I29").Value = 0.062676437182814557
$ws.Range("J29").Value = 0.098516366541463785
$ws.Range("K29").Value = 1
$ws.Range("A30").Value = "X1 on 2"
$ws.Range("B30").Value = -0.0053851758744829503
$ws.Range("C30").Value = -0.28363375249536382
$ws.Range("D30").Value = -0.91114256184009346
$ws.Range("E30").Value = 0.34387505684936581
$ws.Range("F30").Value = 0.32016343886644816
$ws.Range("G30").Value = 0.030253108311141499
$ws.Range("H30").Value = 0.094492701659670678
$ws.Range("I30").Value = 0.094492701659670678
$ws.Range("J30").Value = 0.31853731615374709
$ws.Range("K30").Value = 1
$ws.Range("A31").Value = "X1 on 3"
$ws.Range("B31").Value = -0.52914370863395443
$ws.Range("C31").Value = -0.32802310080378605
$ws.Range("D31").Value = -0.87847023830571236
$ws.Range("E31").Value = 0.2224240366981402
$ws.Range("F31").Value = 0.28084553687913821
$ws.Range("G31").Value = 0.029768115449742335
$ws.Range("H31").Value = 0.10599461818242471
$ws.Range("I31").Value = 0.10599461818242471
$ws.Range("J31").Value = 0.22190530048377166
$ws.Range("K31").Value = 1
$ws.Range("A32").Value = "X2 on 1"
$ws.Range("B32").Value = -0.0000000000000045253302270781528
$ws.Range("C32").Value = -0.00000000000035239319302346658
$ws.Range("D32").Value = -0.000000000010995990710866004
$ws.Range("E32").Value = 0.000000000010291204324819071
$ws.Range("F32").Value = 0.0000000000054305066836931086
$ws.Range("G32").Value = 0.0000000000011225254384841388
$ws.Range("H32").Value = 0.20670731183425159
$ws.Range("I32").Value = 0.20670731183425159
$ws.Range("J32").Value = 0.00000000000037278295569788652
$ws.Range("K32").Value = 1
$ws.Range("A33").Value = "X2 on 2"
$ws.Range("B33").Value = -0.35873501092252547
$ws.Range("C33").Value = -0.29420221881195269
$ws.Range("D33").Value = -0.61957323529194808
$ws.Range("E33").Value = 0.031168797668042633
$ws.Range("F33").Value = 0.16600867110134701
$ws.Range("G33").Value = 0.049534897823508311
$ws.Range("H33").Value = 0.29838741250610723
$ws.Range("I33").Value = 0.29838741250610723
$ws.Range("J33").Value = 0.071057341958787198
$ws.Range("K33").Value = 1
$ws.Range("A34").Value = "X2 on 3"
$ws.Range("B34").Value = 0.35873501092253002
$ws.Range("C34").Value = 0.29420221881230507
$ws.Range("D34").Value = -0.031168797600524573
$ws.Range("E34").Value = 0.61957323522513474
$ws.Range("F34").Value = 0.16600867106707817
$ws.Range("G34").Value = 0.049534897800557094
$ws.Range("H34").Value = 0.29838741242944961
$ws.Range("I34").Value = 0.29838741242944961
$ws.Range("J34").Value = 0.071057341958527406
$ws.Range("K34").Value = 1
$ws.Range("A35").Value = "X3 on 1"
$ws.Range("B35").Value = -0.0000000000000022626651135390764
$ws.Range("C35").Value = -0.00000000000016698456474634105
$ws.Range("D35").Value = -0.0000000000051809124542047919
$ws.Range("E35").Value = 0.0000000000048469433247121091
$ws.Range("F35").Value = 0.0000000000025581734812515308
$ws.Range("G35").Value = 0.00000000000058673350962629082
$ws.Range("H35").Value = 0.22935641930712386
$ws.Range("I35").Value = 0.22935641930712386
$ws.Range("J35").Value = 0.00000000000017777337971269686
$ws.Range("K35").Value = 1
$ws.Range("A36").Value = "X3 on 2"
$ws.Range("B36").Value = -0.17936750546126273
$ws.Range("C36").Value = -0.13849399394658732
$ws.Range("D36").Value = -0.29994295319468273
$ws.Range("E36").Value = 0.02295496530150809
$ws.Range("F36").Value = 0.082373431614858356
$ws.Range("G36").Value = 0.026137110720819298
$ws.Range("H36").Value = 0.31730025335140633
$ws.Range("I36").Value = 0.31730025335140633
$ws.Range("J36").Value = 0.044364599024101371
$ws.Range("K36").Value = 1
$ws.Range("A37").Value = "X3 on 3"
$ws.Range("B37").Value = 0.17936750546126501
$ws.Range("C37").Value = 0.13849399394675432
$ws.Range("D37").Value = -0.022954965669756598
$ws.Range("E37").Value = 0.29994295356326528
$ws.Range("F37").Value = 0.082373431802828911
$ws.Range("G37").Value = 0.026137110933428624
$ws.Range("H37").Value = 0.31730025520839122
$ws.Range("I37").Value = 0.31730025520839122
$ws.Range("J37").Value = 0.044364599023975618
$ws.Range("K37").Value = 1
$ws.Range("A38").Value = "X4 on 1"
$ws.Range("B38").Value = -0.18417155354769044
$ws.Range("C38").Value = -0.27800868717520277
$ws.Range("D38").Value = -0.46296786557083358
$ws.Range("E38").Value = -0.093049508779571993
$ws.Range("F38").Value = 0.094368661799178577
$ws.Range("G38").Value = 0.027310280657375997
$ws.Range("H38").Value = 0.28939989331939131
$ws.Range("I38").Value = 0.28939989331939131
$ws.Range("J38").Value = 0.10918088361686869
$ws.Range("K38").Value = 1
$ws.Range("A39").Value = "X4 on 2"
$ws.Range("B39").Value = 0.18417155354769033
$ws.Range("C39").Value = 0.27800868717520277
$ws.Range("D39").Value = 0.093049508779571993
$ws.Range("E39").Value = 0.46296786557083358
$ws.Range("F39").Value = 0.094368661799178577
$ws.Range("G39").Value = 0.027310280657375997
$ws.Range("H39").Value = 0.28939989331939131
$ws.Range("I39").Value = 0.28939989331939131
$ws.Range("J39").Value = 0.10918088361686878
$ws.Range("K39").Value = 1
$ws.Range("A40").Value = "X4 on 3"
$ws.Range("B40").Value = 0.0000000000000001142046133388892
$ws.Range("C40").Value = 0.0000000000000000044151415997608071
$ws.Range("D40").Value = -0.0000000000000001454347951899156
$ws.Range("E40").Value = 0.00000000000000015426507838943721
$ws.Range("F40").Value = 0.000000000000000076455454269401705
$ws.Range("G40").Value = 0.000000000000000066573622366394219
$ws.Range("H40").Value = 0.8707504651246013
$ws.Range("I40").Value = 0.8707504651246013
$ws.Range("J40").Value = 0.00000000000000010987821249047877
$ws.Range("K40").Value = 0.5
$ws.Range("A41").Value = "X5 on 1"
$ws.Range("B41").Value = 0.36834310709538087
$ws.Range("C41").Value = 0.53352675070691258
$ws.Range("D41").Value = 0.1976388733567804
$ws.Range("E41").Value = 0.86941462805704472
$ws.Range("F41").Value = 0.1713745150418951
$ws.Range("G41").Value = 0.03538979099640506
$ws.Range("H41").Value = 0.20650556465617709
$ws.Range("I41").Value = 0.20650556465617709
$ws.Range("J41").Value = 0.1813052456166325
$ws.Range("K41").Value = 1
$ws.Range("A42").Value = "X5 on 2"
$ws.Range("B42").Value = -0.36834310709538065
$ws.Range("C42").Value = -0.53352675070691258
$ws.Range("D42").Value = -0.86941462805704472
$ws.Range("E42").Value = -0.1976388733567804
$ws.Range("F42").Value = 0.1713745150418951
$ws.Range("G42").Value = 0.03538979099640506
$ws.Range("H42").Value = 0.20650556465617709
$ws.Range("I42").Value = 0.20650556465617709
$ws.Range("J42").Value = 0.18130524561663275
$ws.Range("K42").Value = 1
$ws.Range("A43").Value = "X5 on 3"
$ws.Range("B43").Value = -0.0000000000000002284092266777784
$ws.Range("C43").Value = -0.0000000000000000091163292786210058
$ws.Range("D43").Value = -0.00000000000000030087638378620648
$ws.Range("E43").Value = 0.00000000000000028264372522896452
$ws.Range("F43").Value = 0.00000000000000014885990600284072
$ws.Range("G43").Value = 0.00000000000000014635811821655043
$ws.Range("H43").Value = 0.9831936761652762
$ws.Range("I43").Value = 0.9831936761652762
$ws.Range("J43").Value = 0.00000000000000021948230523036126
$ws.Range("K43").Value = 0.5
$ws.Range("B5").ClearFormats()
$ws.Range("B19").ClearFormats()
$ws.Range("B28").ClearFormats()
